# Applies the "stuff at the bottom of the sheets" commit:
#   1. Adds a new "generic" carrier-kind column (J) to the practice rows 2-5.
#   2. Appends a new "stim details" block (rows 27-36) describing which
#      stimuli still need audio/video/images recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New column J ("generic") for the practice pair rows 2-5 ---
foreach ($row in 2..5) {
    $ws.Cells.Item($row, 10).Value = "generic"   # column J
}

# --- 2. New "stim details" block starting at row 27 ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# month/need counts for the "video" and "audio" word types
$stimRows = @(
    @{ Row = 29; Month = 6; WordType = "video" },
    @{ Row = 30; Month = 6; WordType = "video" },
    @{ Row = 31; Month = 7; WordType = "video" },
    @{ Row = 32; Month = 7; WordType = "video" },
    @{ Row = 33; Month = 6; WordType = "audio" },
    @{ Row = 34; Month = 6; WordType = "audio" },
    @{ Row = 35; Month = 7; WordType = "audio" },
    @{ Row = 36; Month = 7; WordType = "audio" }
)

foreach ($entry in $stimRows) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Month     # column A
    $ws.Cells.Item($entry.Row, 2).Value = $entry.WordType  # column B
}
